$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before: ... / 2022-6-26 星期日 / 多云，终于到了周末，我要出去嗨！(+bookmark) / (empty)
# After:  ... / 2022-6-26 星期日 / 多云，终于到了周末，我要出去嗨！/ 2022-7-1 星期五 /
#              今天天气真不错，心情也好 (+bookmark)
#
# i.e. two new diary-entry paragraphs are inserted right before the final
# (bookmarked) paragraph, that final paragraph's text is changed, and the
# trailing empty paragraph is removed.
# ---------------------------------------------------------------------------

# Paragraph 5 ("2022-6-26星期日") already carries the eastAsia font hint we
# want the freshly inserted paragraphs to use, so insert after it - Word
# copies the anchor paragraph's formatting onto new paragraphs it creates.
$anchor = $d.Paragraphs.Item(5).Range
$anchor.InsertParagraphAfter()
$d.Paragraphs.Item(6).Range.Text = "多云，终于到了周末，我要出去嗨！"

$anchor2 = $d.Paragraphs.Item(6).Range
$anchor2.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.Text = "2022-7-1 星期五"

# Paragraph 8 is now the original last-content paragraph (with the
# _GoBack bookmark); update its text in place so the bookmark survives.
$targetRange = $d.Paragraphs.Item(8).Range
[void]$targetRange.Find.Execute("多云，终于到了周末，我要出去嗨！", $false, $false, $false, $false, $false,
                           $true, 1, $false, "今天天气真不错，心情也好", 2)

# Remove the now-redundant trailing empty paragraph (paragraph 9), including
# the paragraph mark that precedes it so the paragraph count actually drops.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$d.Range($lastPara.Range.Start - 1, $lastPara.Range.End).Delete()
